# Remove the "url" field row from the Documents[]/File group (row 12)
# on the "Specification" sheet. Deleting this entire row shifts every
# subsequent row up by one, which matches the new sheet dimension
# (A1:N111 -> A1:N110) and the renumbered merged cell ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Specification")

$ws.Rows.Item(12).Delete()
